# Update the two-digit / one-digit division answers in the answer table.
# The table has 5 data rows (1, 5, 9, 13, 17) x 5 columns; the other rows
# are blank spacer rows. Each populated cell gets its old "a÷b=c, d" text
# replaced with a freshly generated one, keeping the existing run
# formatting (TimeNewRoman, sz 30) intact by only touching the cell's
# Range.Text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "69÷2=34, 1" },
    @{ Row = 1;  Col = 2; New = "40÷8=5, 0" },
    @{ Row = 1;  Col = 3; New = "74÷9=8, 2" },
    @{ Row = 1;  Col = 4; New = "80÷9=8, 8" },
    @{ Row = 1;  Col = 5; New = "10÷5=2, 0" },

    @{ Row = 5;  Col = 1; New = "90÷2=45, 0" },
    @{ Row = 5;  Col = 2; New = "50÷4=12, 2" },
    @{ Row = 5;  Col = 3; New = "57÷2=28, 1" },
    @{ Row = 5;  Col = 4; New = "77÷3=25, 2" },
    @{ Row = 5;  Col = 5; New = "77÷7=11, 0" },

    @{ Row = 9;  Col = 1; New = "41÷7=5, 6" },
    @{ Row = 9;  Col = 2; New = "99÷8=12, 3" },
    @{ Row = 9;  Col = 3; New = "35÷4=8, 3" },
    @{ Row = 9;  Col = 4; New = "76÷6=12, 4" },
    @{ Row = 9;  Col = 5; New = "18÷3=6, 0" },

    @{ Row = 13; Col = 1; New = "50÷9=5, 5" },
    @{ Row = 13; Col = 2; New = "59÷5=11, 4" },
    @{ Row = 13; Col = 3; New = "18÷4=4, 2" },
    @{ Row = 13; Col = 4; New = "55÷4=13, 3" },
    @{ Row = 13; Col = 5; New = "40÷5=8, 0" },

    @{ Row = 17; Col = 1; New = "33÷8=4, 1" },
    @{ Row = 17; Col = 2; New = "53÷9=5, 8" },
    @{ Row = 17; Col = 3; New = "11÷7=1, 4" },
    @{ Row = 17; Col = 4; New = "93÷9=10, 3" },
    @{ Row = 17; Col = 5; New = "27÷5=5, 2" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
